# Fix: swap the shared-string values in column B so that
# B1 -> "hi" and B2 -> "hello" (previously B1 -> "hello", B2 -> "hi").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "hi"
$ws.Range("B2").Value = "hello"

# Restore the default selection to A1 (removes the stale D2 selection
# left over from the previous session).
$ws.Range("A1").Select() | Out-Null
